$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 198, pushing all existing rows (198..245) down to (199..246).
$ws.Rows.Item(198).Insert()

# Populate the newly inserted row 198 with the new weekly price entry.
$ws.Cells.Item(198, 1).Value = 4
$ws.Cells.Item(198, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(198, 3).Value = "Los Lagos"
$ws.Cells.Item(198, 4).Value = 44543
$ws.Cells.Item(198, 5).Value = 10
$ws.Cells.Item(198, 6).Value = 100112008
$ws.Cells.Item(198, 7).Value = "Coliflor"
$ws.Cells.Item(198, 8).Value = "Sin especificar"
$ws.Cells.Item(198, 9).Value = "Primera"
$ws.Cells.Item(198, 10).Value = 500
$ws.Cells.Item(198, 11).Value = 1300
$ws.Cells.Item(198, 12).Value = 1300
$ws.Cells.Item(198, 13).Value = 1300
$ws.Cells.Item(198, 14).Value = "$/unidad"
$ws.Cells.Item(198, 15).Value = "Región Metropolitana"
$ws.Cells.Item(198, 16).Value = 1300
$ws.Cells.Item(198, 17).Value = 1
$ws.Cells.Item(198, 18).Value = "Hortaliza"
